$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.115.62"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "1.824.06"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5989"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06948"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2749"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07595"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("D12").Value = "1.832.33"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.735"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6257"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009785"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.08%  "
$ws.Range("D17").Value = "28.781.12"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.578"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -10.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.879"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.48%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "156.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.927"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1286"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.52%  "
$ws.Range("E27").Value = "  -4.60%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.427"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.66%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06379"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -10.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.440"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.828"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.755"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.728"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.089"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6471"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.531"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("E38").Value = "  -4.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.542"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.41%  "
$ws.Range("D40").Value = "1.144.54"
$ws.Range("E40").Value = "  -7.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8850"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.36%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "1.980.39"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000112"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.600"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.488"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05498"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4538"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.415"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.20%  "
